$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated voltage magnitude (vm_pu) results for the 380 kV case
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.0241563437844
$ws.Cells.Item(2, 4).Value = 1.044899987921914
$ws.Cells.Item(2, 5).Value = 1.035976766366692
$ws.Cells.Item(2, 6).Value = 1.05014729415626
$ws.Cells.Item(2, 9).Value = 1.0370446962429
$ws.Cells.Item(2, 10).Value = 1.029332536153319
$ws.Cells.Item(2, 11).Value = 1.047669739409689
$ws.Cells.Item(2, 12).Value = 1.038771828197156
$ws.Cells.Item(2, 13).Value = 1.052902379673941
$ws.Cells.Item(2, 14).Value = 1.013799798924982

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.02503909946454
$ws.Cells.Item(3, 4).Value = 1.04542917025987
$ws.Cells.Item(3, 5).Value = 1.036716946810454
$ws.Cells.Item(3, 6).Value = 1.050867794244345
$ws.Cells.Item(3, 9).Value = 1.037123224886164
$ws.Cells.Item(3, 10).Value = 1.029854423236858
$ws.Cells.Item(3, 11).Value = 1.048010989095989
$ws.Cells.Item(3, 12).Value = 1.039321698991654
$ws.Cells.Item(3, 13).Value = 1.053435504954948
$ws.Cells.Item(3, 14).Value = 1.013975598377531

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.025611230685561
$ws.Cells.Item(4, 4).Value = 1.045771660964924
$ws.Cells.Item(4, 5).Value = 1.037196926358242
$ws.Cells.Item(4, 6).Value = 1.051334660772645
$ws.Cells.Item(4, 9).Value = 1.037172754600657
$ws.Cells.Item(4, 10).Value = 1.030192423861757
$ws.Cells.Item(4, 11).Value = 1.048231162733385
$ws.Cells.Item(4, 12).Value = 1.039677867915766
$ws.Cells.Item(4, 13).Value = 1.053780454368193
$ws.Cells.Item(4, 14).Value = 1.014089375676072

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.025851975595366
$ws.Cells.Item(5, 4).Value = 1.045915659738252
$ws.Cells.Item(5, 5).Value = 1.037398955232417
$ws.Cells.Item(5, 6).Value = 1.051531086035978
$ws.Cells.Item(5, 9).Value = 1.037193268828092
$ws.Cells.Item(5, 10).Value = 1.030334591028399
$ws.Cells.Item(5, 11).Value = 1.048323569500357
$ws.Cells.Item(5, 12).Value = 1.039827687405604
$ws.Cells.Item(5, 13).Value = 1.05392546494729
$ws.Cells.Item(5, 14).Value = 1.014137212738823

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.025892410637137
$ws.Cells.Item(6, 4).Value = 1.045939838623742
$ws.Cells.Item(6, 5).Value = 1.037432891108605
$ws.Cells.Item(6, 6).Value = 1.051564075702778
$ws.Cells.Item(6, 9).Value = 1.03719669516953
$ws.Cells.Item(6, 10).Value = 1.030358465689318
$ws.Cells.Item(6, 11).Value = 1.04833907591367
$ws.Cells.Item(6, 12).Value = 1.039852847761352
$ws.Cells.Item(6, 13).Value = 1.053949812449053
$ws.Cells.Item(6, 14).Value = 1.014145245076145

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.025614446664568
$ws.Cells.Item(7, 4).Value = 1.045773585024645
$ws.Cells.Item(7, 5).Value = 1.037199624915805
$ws.Cells.Item(7, 6).Value = 1.051337284812009
$ws.Cells.Item(7, 9).Value = 1.037173029924225
$ws.Cells.Item(7, 10).Value = 1.030194323226232
$ws.Cells.Item(7, 11).Value = 1.048232398085114
$ws.Cells.Item(7, 12).Value = 1.039679869474781
$ws.Cells.Item(7, 13).Value = 1.053782392033069
$ws.Cells.Item(7, 14).Value = 1.014090014857876

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.024454481870679
$ws.Cells.Item(8, 4).Value = 1.045078810905115
$ws.Cells.Item(8, 5).Value = 1.036226698716353
$ws.Cells.Item(8, 6).Value = 1.050390653643122
$ws.Cells.Item(8, 9).Value = 1.037071500670842
$ws.Cells.Item(8, 10).Value = 1.029508846470091
$ws.Cells.Item(8, 11).Value = 1.047785197449867
$ws.Cells.Item(8, 12).Value = 1.038957582981297
$ws.Cells.Item(8, 13).Value = 1.053082554579882
$ws.Cells.Item(8, 14).Value = 1.013859206054033

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.022417660345537
$ws.Cells.Item(9, 4).Value = 1.043855193017705
$ws.Cells.Item(9, 5).Value = 1.034520276494633
$ws.Cells.Item(9, 6).Value = 1.048727680012546
$ws.Cells.Item(9, 9).Value = 1.036882796826276
$ws.Cells.Item(9, 10).Value = 1.028303338968049
$ws.Cells.Item(9, 11).Value = 1.046992357286245
$ws.Cells.Item(9, 12).Value = 1.037687688888442
$ws.Cells.Item(9, 13).Value = 1.051849289906287
$ws.Cells.Item(9, 14).Value = 1.013452691591164

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.021064700782483
$ws.Cells.Item(10, 4).Value = 1.043040019830029
$ws.Cells.Item(10, 5).Value = 1.033388152813631
$ws.Cells.Item(10, 6).Value = 1.047622598003461
$ws.Cells.Item(10, 9).Value = 1.036750454067987
$ws.Cells.Item(10, 10).Value = 1.027501349598762
$ws.Cells.Item(10, 11).Value = 1.046460652385442
$ws.Cells.Item(10, 12).Value = 1.03684310923902
$ws.Cells.Item(10, 13).Value = 1.051027177971444
$ws.Cells.Item(10, 14).Value = 1.013181845804156

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.020480040763403
$ws.Cells.Item(11, 4).Value = 1.042687203315955
$ws.Cells.Item(11, 5).Value = 1.032899257251636
$ws.Cells.Item(11, 6).Value = 1.047144958983862
$ws.Cells.Item(11, 9).Value = 1.036691606725691
$ws.Cells.Item(11, 10).Value = 1.02715449426018
$ws.Cells.Item(11, 11).Value = 1.046229692050221
$ws.Cells.Item(11, 12).Value = 1.036477893967872
$ws.Cells.Item(11, 13).Value = 1.050671232156122
$ws.Cells.Item(11, 14).Value = 1.013064611637855

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.020263050874145
$ws.Cells.Item(12, 4).Value = 1.042556177598031
$ws.Cells.Item(12, 5).Value = 1.032717860165144
$ws.Cells.Item(12, 6).Value = 1.04696767535641
$ws.Cells.Item(12, 9).Value = 1.036669517301712
$ws.Cells.Item(12, 10).Value = 1.027025719826033
$ws.Cells.Item(12, 11).Value = 1.04614379522042
$ws.Cells.Item(12, 12).Value = 1.036342312386045
$ws.Cells.Item(12, 13).Value = 1.050539024966738
$ws.Cells.Item(12, 14).Value = 1.013021072734841

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.020309587822909
$ws.Cells.Item(13, 4).Value = 1.042584281845833
$ws.Cells.Item(13, 5).Value = 1.032756761351795
$ws.Cells.Item(13, 6).Value = 1.047005697252093
$ws.Cells.Item(13, 9).Value = 1.036674265994449
$ws.Cells.Item(13, 10).Value = 1.027053339504102
$ws.Cells.Item(13, 11).Value = 1.046162225238259
$ws.Cells.Item(13, 12).Value = 1.036371391648732
$ws.Cells.Item(13, 13).Value = 1.050567383514715
$ws.Cells.Item(13, 14).Value = 1.013030411648663

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.020462100653615
$ws.Cells.Item(14, 4).Value = 1.042676372145367
$ws.Cells.Item(14, 5).Value = 1.032884258810558
$ws.Cells.Item(14, 6).Value = 1.047130301942442
$ws.Cells.Item(14, 9).Value = 1.03668978551284
$ws.Cells.Item(14, 10).Value = 1.027143848431008
$ws.Cells.Item(14, 11).Value = 1.046222593982987
$ws.Cells.Item(14, 12).Value = 1.036466685200386
$ws.Cells.Item(14, 13).Value = 1.050660303716648
$ws.Cells.Item(14, 14).Value = 1.01306101255111

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.020556092590609
$ws.Cells.Item(15, 4).Value = 1.042733115534456
$ws.Cells.Item(15, 5).Value = 1.032962840828702
$ws.Cells.Item(15, 6).Value = 1.04720709266525
$ws.Cells.Item(15, 9).Value = 1.036699317026905
$ws.Cells.Item(15, 10).Value = 1.027199622365411
$ws.Cells.Item(15, 11).Value = 1.046259774911569
$ws.Cells.Item(15, 12).Value = 1.036525408769939
$ws.Cells.Item(15, 13).Value = 1.050717555892998
$ws.Cells.Item(15, 14).Value = 1.013079867730967

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.02110352774274
$ws.Cells.Item(16, 4).Value = 1.043063438633106
$ws.Cells.Item(16, 5).Value = 1.033420627210863
$ws.Cells.Item(16, 6).Value = 1.047654315828233
$ws.Cells.Item(16, 9).Value = 1.03675432714917
$ws.Cells.Item(16, 10).Value = 1.027524378021307
$ws.Cells.Item(16, 11).Value = 1.046475965246924
$ws.Cells.Item(16, 12).Value = 1.03686735789657
$ws.Cells.Item(16, 13).Value = 1.051050801795908
$ws.Cells.Item(16, 14).Value = 1.01318962721605

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.021447236457059
$ws.Cells.Item(17, 4).Value = 1.043270685656049
$ws.Cells.Item(17, 5).Value = 1.033708139690123
$ws.Cells.Item(17, 6).Value = 1.04793508128787
$ws.Cells.Item(17, 9).Value = 1.036788421241301
$ws.Cells.Item(17, 10).Value = 1.02772819969908
$ws.Cells.Item(17, 11).Value = 1.046611381785989
$ws.Cells.Item(17, 12).Value = 1.037081986690549
$ws.Cells.Item(17, 13).Value = 1.051259848552376
$ws.Cells.Item(17, 14).Value = 1.013258488532406

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.021647829723463
$ws.Cells.Item(18, 4).Value = 1.043391584519742
$ws.Cells.Item(18, 5).Value = 1.033875968225704
$ws.Cells.Item(18, 6).Value = 1.048098930704784
$ws.Cells.Item(18, 9).Value = 1.036808158910536
$ws.Cells.Item(18, 10).Value = 1.027847125038726
$ws.Cells.Item(18, 11).Value = 1.046690297556189
$ws.Cells.Item(18, 12).Value = 1.037207223528341
$ws.Cells.Item(18, 13).Value = 1.051381785155034
$ws.Cells.Item(18, 14).Value = 1.013298658370342

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.021716246097475
$ws.Cells.Item(19, 4).Value = 1.043432810416635
$ws.Cells.Item(19, 5).Value = 1.033933214963939
$ws.Cells.Item(19, 6).Value = 1.048154813245948
$ws.Cells.Item(19, 9).Value = 1.036814863676592
$ws.Cells.Item(19, 10).Value = 1.027887682184684
$ws.Cells.Item(19, 11).Value = 1.046717193804056
$ws.Cells.Item(19, 12).Value = 1.037249934063091
$ws.Cells.Item(19, 13).Value = 1.051423362867233
$ws.Cells.Item(19, 14).Value = 1.013312355943729

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.021410347984114
$ws.Cells.Item(20, 4).Value = 1.043248448425834
$ws.Cells.Item(20, 5).Value = 1.033677279133406
$ws.Cells.Item(20, 6).Value = 1.047904949140135
$ws.Cells.Item(20, 9).Value = 1.036784778657285
$ws.Cells.Item(20, 10).Value = 1.027706327456393
$ws.Cells.Item(20, 11).Value = 1.046596860143884
$ws.Cells.Item(20, 12).Value = 1.037058954128741
$ws.Cells.Item(20, 13).Value = 1.051237419473191
$ws.Cells.Item(20, 14).Value = 1.013251099927129

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.020417184476132
$ws.Cells.Item(21, 4).Value = 1.042649253111031
$ws.Cells.Item(21, 5).Value = 1.032846708436429
$ws.Cells.Item(21, 6).Value = 1.047093605283666
$ws.Cells.Item(21, 9).Value = 1.036685221769258
$ws.Cells.Item(21, 10).Value = 1.027117194056234
$ws.Cells.Item(21, 11).Value = 1.046204819857549
$ws.Cells.Item(21, 12).Value = 1.036438621524566
$ws.Cells.Item(21, 13).Value = 1.050632940827155
$ws.Cells.Item(21, 14).Value = 1.013052001148429

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.019793778831186
$ws.Cells.Item(22, 4).Value = 1.042272666993274
$ws.Cells.Item(22, 5).Value = 1.032325656634368
$ws.Cells.Item(22, 6).Value = 1.046584250481816
$ws.Cells.Item(22, 9).Value = 1.036621290864428
$ws.Cells.Item(22, 10).Value = 1.026747148114643
$ws.Cells.Item(22, 11).Value = 1.045957705204288
$ws.Cells.Item(22, 12).Value = 1.036049032125113
$ws.Cells.Item(22, 13).Value = 1.050252922130349
$ws.Cells.Item(22, 14).Value = 1.012926861089313

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.02012415927614
$ws.Cells.Item(23, 4).Value = 1.042472287309159
$ws.Cells.Item(23, 5).Value = 1.032601765354221
$ws.Cells.Item(23, 6).Value = 1.046854195431403
$ws.Cells.Item(23, 9).Value = 1.036655308181523
$ws.Cells.Item(23, 10).Value = 1.026943281386127
$ws.Cells.Item(23, 11).Value = 1.046088763916618
$ws.Cells.Item(23, 12).Value = 1.036255518781231
$ws.Cells.Item(23, 13).Value = 1.050454372738201
$ws.Cells.Item(23, 14).Value = 1.012993196126344

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.021427015945783
$ws.Cells.Item(24, 4).Value = 1.043258496426353
$ws.Cells.Item(24, 5).Value = 1.033691223297309
$ws.Cells.Item(24, 6).Value = 1.047918564302422
$ws.Cells.Item(24, 9).Value = 1.036786425044189
$ws.Cells.Item(24, 10).Value = 1.027716210459174
$ws.Cells.Item(24, 11).Value = 1.046603422066618
$ws.Cells.Item(24, 12).Value = 1.037069361405442
$ws.Cells.Item(24, 13).Value = 1.051247554199071
$ws.Cells.Item(24, 14).Value = 1.01325443850691

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.022943366767435
$ws.Cells.Item(25, 4).Value = 1.044171435968599
$ws.Cells.Item(25, 5).Value = 1.034960468010078
$ws.Cells.Item(25, 6).Value = 1.049156979372135
$ws.Cells.Item(25, 9).Value = 1.036932737435465
$ws.Cells.Item(25, 10).Value = 1.028614700553285
$ws.Cells.Item(25, 11).Value = 1.047197886550329
$ws.Cells.Item(25, 12).Value = 1.038015638043611
$ws.Cells.Item(25, 13).Value = 1.052168114764148
$ws.Cells.Item(25, 14).Value = 1.013557758587352
